$wb = $excel.ActiveWorkbook

# ALC sheet updates
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 114210.78
$ws.Range("J32").Value = 202619.4
$ws.Range("L32").Value = 202619.4
$ws.Range("N32").Value = -203271.4
$ws.Range("H33").Value = 358
$ws.Range("I33").Value = 340.25
$ws.Range("K33").Value = 340.25
$ws.Range("M33").Value = -111.25
$ws.Range("H92").Value = 348.26666
$ws.Range("I92").Value = 325
$ws.Range("K92").Value = 325
$ws.Range("M92").Value = 923
$ws.Range("H112").Value = 1297.7391
$ws.Range("J112").Value = 1326.2632
$ws.Range("L112").Value = 3978.7896
$ws.Range("N112").Value = -6194.7896
$ws.Range("H127").Value = 587.3333
$ws.Range("I127").Value = 472.5
$ws.Range("J127").Value = 817
$ws.Range("K127").Value = 1417.5
$ws.Range("L127").Value = 2451
$ws.Range("M127").Value = 3542.5
$ws.Range("N127").Value = -12371
$ws.Range("H132").Value = 1973.4474
$ws.Range("I132").Value = 1527.6389
$ws.Range("K132").Value = 4582.9167
$ws.Range("M132").Value = -2052.9167
$ws.Range("H137").Value = 693779.2
$ws.Range("I137").Value = 2360.818
$ws.Range("K137").Value = 7082.454000000001
$ws.Range("M137").Value = -4532.454000000001
$ws.Range("H138").Value = 1791.25
$ws.Range("I138").Value = 1455.8948
$ws.Range("K138").Value = 4367.6844
$ws.Range("M138").Value = 772.3155999999999
$ws.Range("H141").Value = 2605.76
$ws.Range("I141").Value = 2267.1304
$ws.Range("K141").Value = 6801.3912
$ws.Range("M141").Value = -1621.3912


# ARM sheet updates
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5889.0728
$ws.Range("I32").Value = 2272.6
$ws.Range("K32").Value = 2272.6
$ws.Range("M32").Value = -1985.6
$ws.Range("H45").Value = 10419027
$ws.Range("J45").Value = 31251996
$ws.Range("L45").Value = 31251996
$ws.Range("N45").Value = -31252750
$ws.Range("H74").Value = 146378.86
$ws.Range("I74").Value = 251225.5
$ws.Range("J74").Value = 6583.3335
$ws.Range("K74").Value = 251225.5
$ws.Range("L74").Value = 6583.3335
$ws.Range("M74").Value = -250351.5
$ws.Range("N74").Value = -8331.333500000001
$ws.Range("H77").Value = 146378.86
$ws.Range("I77").Value = 251225.5
$ws.Range("J77").Value = 6583.3335
$ws.Range("K77").Value = 1256127.5
$ws.Range("L77").Value = 32916.6675
$ws.Range("M77").Value = -1251759.5
$ws.Range("N77").Value = -41652.6675
$ws.Range("H102").Value = 48040.543
$ws.Range("I102").Value = 51860.65
$ws.Range("K102").Value = 51860.65
$ws.Range("M102").Value = -50238.65
$ws.Range("H120").Value = 139990
$ws.Range("J120").Value = 139990
$ws.Range("L120").Value = 139990
$ws.Range("N120").Value = -149666


# BSM sheet updates
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1277.3529
$ws.Range("I20").Value = 1167.875
$ws.Range("J20").Value = 1374.6666
$ws.Range("K20").Value = 1167.875
$ws.Range("L20").Value = 1374.6666
$ws.Range("M20").Value = -920.875
$ws.Range("N20").Value = -1868.6666
$ws.Range("H86").Value = 4738.0586
$ws.Range("I86").Value = 4257.5557
$ws.Range("K86").Value = 4257.5557
$ws.Range("M86").Value = -3134.5557
$ws.Range("H89").Value = 4738.0586
$ws.Range("I89").Value = 4257.5557
$ws.Range("K89").Value = 21287.7785
$ws.Range("M89").Value = -15671.7785
$ws.Range("H112").Value = 91658.336
$ws.Range("J112").Value = 91658.336
$ws.Range("L112").Value = 91658.336
$ws.Range("N112").Value = -94612.336
$ws.Range("H126").Value = 65000
$ws.Range("J126").Value = 65000
$ws.Range("L126").Value = 65000
$ws.Range("N126").Value = -74880
$ws.Range("H134").Value = 7281.1665
$ws.Range("I134").Value = 3422
$ws.Range("K134").Value = 10266
$ws.Range("M134").Value = -7731


# CRP sheet updates
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2374.9033
$ws.Range("I31").Value = 1413.4348
$ws.Range("J31").Value = 5139.125
$ws.Range("K31").Value = 1413.4348
$ws.Range("L31").Value = 5139.125
$ws.Range("M31").Value = -1118.4348
$ws.Range("N31").Value = -5729.125
$ws.Range("H34").Value = 2374.9033
$ws.Range("I34").Value = 1413.4348
$ws.Range("J34").Value = 5139.125
$ws.Range("K34").Value = 1413.4348
$ws.Range("L34").Value = 5139.125
$ws.Range("M34").Value = -1211.4348
$ws.Range("N34").Value = -5543.125
$ws.Range("H58").Value = 1634.8182
$ws.Range("I58").Value = 1475.2778
$ws.Range("J58").Value = 2352.75
$ws.Range("K58").Value = 1475.2778
$ws.Range("L58").Value = 2352.75
$ws.Range("M58").Value = -1272.2778
$ws.Range("N58").Value = -2758.75
$ws.Range("H76").Value = 6666.6665
$ws.Range("I76").Value = 6666.6665
$ws.Range("K76").Value = 6666.6665
$ws.Range("M76").Value = -6351.6665
$ws.Range("H79").Value = 6666.6665
$ws.Range("I79").Value = 6666.6665
$ws.Range("K79").Value = 6666.6665
$ws.Range("M79").Value = -5574.6665
$ws.Range("H99").Value = 6946740.5
$ws.Range("I99").Value = 11113321
$ws.Range("J99").Value = 2440.1667
$ws.Range("K99").Value = 11113321
$ws.Range("L99").Value = 2440.1667
$ws.Range("M99").Value = -11111823
$ws.Range("N99").Value = -5436.1667
$ws.Range("H105").Value = 4381.231
$ws.Range("I105").Value = 3117.4443
$ws.Range("J105").Value = 7224.75
$ws.Range("K105").Value = 3117.4443
$ws.Range("L105").Value = 7224.75
$ws.Range("M105").Value = -1370.4443
$ws.Range("N105").Value = -10718.75
$ws.Range("H122").Value = 3348.875
$ws.Range("I122").Value = 2505.4285
$ws.Range("K122").Value = 7516.2855
$ws.Range("M122").Value = -5066.2855
$ws.Range("H126").Value = 6946740.5
$ws.Range("I126").Value = 11113321
$ws.Range("J126").Value = 2440.1667
$ws.Range("K126").Value = 33339963
$ws.Range("L126").Value = 7320.500100000001
$ws.Range("M126").Value = -33337493
$ws.Range("N126").Value = -12260.5001
$ws.Range("H136").Value = 1634.8182
$ws.Range("I136").Value = 1475.2778
$ws.Range("J136").Value = 2352.75
$ws.Range("K136").Value = 4425.8334
$ws.Range("L136").Value = 7058.25
$ws.Range("M136").Value = -1875.8334
$ws.Range("N136").Value = -12158.25


# CUL sheet updates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1547.5
$ws.Range("I52").Value = 790
$ws.Range("J52").Value = 1800
$ws.Range("K52").Value = 2370
$ws.Range("L52").Value = 5400
$ws.Range("M52").Value = -2104
$ws.Range("N52").Value = -5932
$ws.Range("H132").Value = 3289.6155
$ws.Range("J132").Value = 4618
$ws.Range("L132").Value = 41562
$ws.Range("N132").Value = -46622


# GSM sheet updates
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 6947.25
$ws.Range("J18").Value = 9994.5
$ws.Range("L18").Value = 9994.5
$ws.Range("N18").Value = -10580.5


# LTW sheet updates
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2308.3333
$ws.Range("I46").Value = 981.3333
$ws.Range("J46").Value = 2971.8333
$ws.Range("K46").Value = 981.3333
$ws.Range("L46").Value = 2971.8333
$ws.Range("M46").Value = -793.3333
$ws.Range("N46").Value = -3347.8333
$ws.Range("H64").Value = 16000
$ws.Range("J64").Value = 16000
$ws.Range("L64").Value = 16000
$ws.Range("N64").Value = -16450
$ws.Range("H67").Value = 16000
$ws.Range("J67").Value = 16000
$ws.Range("L67").Value = 16000
$ws.Range("N67").Value = -17560
$ws.Range("H132").Value = 3733.3333
$ws.Range("I132").Value = 3733.3333
$ws.Range("K132").Value = 11199.9999
$ws.Range("M132").Value = -8669.999899999999


# WVR sheet updates
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 18249
$ws.Range("J63").Value = 18249
$ws.Range("L63").Value = 18249
$ws.Range("N63").Value = -19497
$ws.Range("H66").Value = 18249
$ws.Range("J66").Value = 18249
$ws.Range("L66").Value = 54747
$ws.Range("N66").Value = -60987
$ws.Range("H132").Value = 3250
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9750
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14810
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

